# Update countries & provincias Spain
# This workbook ("Pais" sheet) lists countries with COVID-19 style stats
# (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes) sorted descending by column B ("Casos totales").
#
# The source data was refreshed and the sheet re-sorted. That re-sort only
# changes row order/content for a handful of rows whose updated totals (or
# ties) moved them relative to their neighbours; all other rows are
# unaffected. Below we just poke the new values into the affected rows
# (same effect as refreshing the data + re-sorting), and bump the
# "last updated" timestamp shown in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 11:05"

# --- Row 31: becomes Banglades with refreshed figures ---
$ws.Cells.Item(31, 1).Value = "Banglades"
$ws.Cells.Item(31, 2).Value = 25121
$ws.Cells.Item(31, 3).Value = 1251
$ws.Cells.Item(31, 4).Value = 4993
$ws.Cells.Item(31, 5).Value = 19758
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 21
$ws.Cells.Item(31, 8).Value = 370

# --- Row 32: becomes Irlanda (previous Irlanda figures) ---
$ws.Cells.Item(32, 1).Value = "Irlanda"
$ws.Cells.Item(32, 2).Value = 24200
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 19470
$ws.Cells.Item(32, 5).Value = 3183
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 1547

# --- Row 33: becomes Emiratos Arabes Unidos (previous Emiratos figures) ---
$ws.Cells.Item(33, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(33, 2).Value = 24190
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 9577
$ws.Cells.Item(33, 5).Value = 14389
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 224

# --- Row 51: Chequia refreshed figures ---
$ws.Cells.Item(51, 1).Value = "Chequia"
$ws.Cells.Item(51, 2).Value = 8594
$ws.Cells.Item(51, 3).Value = 8
$ws.Cells.Item(51, 4).Value = 5641
$ws.Cells.Item(51, 5).Value = 2654
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = 299

# --- Row 58: becomes Malasia with refreshed figures ---
$ws.Cells.Item(58, 1).Value = "Malasia"
$ws.Cells.Item(58, 2).Value = 6978
$ws.Cells.Item(58, 3).Value = 37
$ws.Cells.Item(58, 4).Value = 5646
$ws.Cells.Item(58, 5).Value = 1218
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 114

# --- Row 59: becomes Marruecos (previous Marruecos figures) ---
$ws.Cells.Item(59, 1).Value = "Marruecos"
$ws.Cells.Item(59, 2).Value = 6952
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 3758
$ws.Cells.Item(59, 5).Value = 3002
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 192

# --- Row 196: becomes Nueva Caledonia (tied total, swapped with Belice) ---
$ws.Cells.Item(196, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(196, 2).Value = 18
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 18
$ws.Cells.Item(196, 5).Value = 0
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0

# --- Row 197: becomes Belice (previous Belice figures) ---
$ws.Cells.Item(197, 1).Value = "Belice"
$ws.Cells.Item(197, 2).Value = 18
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 16
$ws.Cells.Item(197, 5).Value = 0
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 2

# --- Row 209: becomes Seychelles (tied total, 3-way rotation) ---
$ws.Cells.Item(209, 1).Value = "Seychelles"
$ws.Cells.Item(209, 2).Value = 11
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# --- Row 210: becomes Groenlandia ---
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# --- Row 211: becomes Montserrat ---
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

# --- Row 215: becomes San Bartolome (tied total, swapped with Bonaire...) ---
$ws.Cells.Item(215, 1).Value = "San Bartolome"
$ws.Cells.Item(215, 2).Value = 6
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 6
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

# --- Row 216: becomes Bonaire, San Eustaquio y Saba ---
$ws.Cells.Item(216, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 2).Value = 6
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 6
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0
